$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the BD / FFC quantity drivers (I2, I3) ---
# BD (Board Qtys) 3 -> 5
$ws.Range("I2").Value = 5
# FFC 3 -> 5
$ws.Range("I3").Value = 5

# --- Row 28: part swap (649-68015-403HLF -> 538-90121-0763) ---
$ws.Range("A28").Value = "538-90121-0763"
$ws.Range("B28").Value = "Headers & Wire Housings 2.54MM CGRIDIII HDR 3P R/A SR SEL AU"
$ws.Range("C28").Value = 0.75
$ws.Range("D28").Formula = "=1*BD+1*FFC"

# --- Row 29: same part, quantity driver formula stays 1*BD (value follows BD) ---
$ws.Range("A29").Value = "649-68021-408HLF"
$ws.Range("B29").Value = "Headers & Wire Housings 8P SR UNSHRD HRD TIN OVER NI"
$ws.Range("C29").Value = 0.47
$ws.Range("D29").Formula = "=1*BD"

# --- Row 30: same part, quantity driver formula becomes 1*BD+1*FFC ---
$ws.Range("A30").Value = "649-77313-418-16LF"
$ws.Range("B30").Value = "Headers & Wire Housings 16P STR DR TMT HDR TIN .45IN LENGTH"
$ws.Range("C30").Value = 0.5
$ws.Range("D30").Formula = "=1*BD+1*FFC"

# --- Row 31: brand-new BOM line ---
$ws.Range("A31").Value = "649-SFW16R-2STE1LF"
$ws.Range("B31").Value = "FFC & FPC Connectors 16P SIDE SMT ZIF UPPER CONTACT"
$ws.Range("C31").Value = 0.73
$ws.Range("D31").Formula = "=FFC*1"
$ws.Range("E31").Formula = "=C31*D31"

# Match formatting/styles with the row above it.
$ws.Range("A30:E30").Copy()
$ws.Range("A31:E31").PasteSpecial(-4122)

# Re-apply the values/formulas for row 31 in case PasteSpecial touched formats only.
$ws.Range("A31").Value = "649-SFW16R-2STE1LF"
$ws.Range("B31").Value = "FFC & FPC Connectors 16P SIDE SMT ZIF UPPER CONTACT"
$ws.Range("C31").Value = 0.73
$ws.Range("D31").Formula = "=FFC*1"
$ws.Range("E31").Formula = "=C31*D31"

# Widen column I slightly (metadata cosmetic change from the diff).
$ws.Columns.Item(9).ColumnWidth = 15

# Move the active selection like the source workbook ended up with.
$ws.Range("D24").Select()

$wb.Application.CalculateFullRebuild()
